$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("EM_Exp1")

# Add new "selex_type" column (M) header
$ws1.Range("M1").Value = "selex_type"

# Existing EM rows (2-14) are length-based selectivity
for ($r = 2; $r -le 14; $r++) {
    $ws1.Cells.Item($r, 13).Value = "length"
}

# Add a new EM row (15), cloned from row 2 ("Age"), representing an
# age-based-selectivity variant of the vonB/age EM.
for ($c = 1; $c -le 12; $c++) {
    $ws1.Cells.Item(15, $c).Value = $ws1.Cells.Item(2, $c).Value2
}
$ws1.Range("M15").Value = "age"
$ws1.Range("A15").Value = "Age_AgeSelex"

# Make EM_Exp1 the active sheet/tab, with A16 selected (just below new data)
$ws1.Activate() | Out-Null
$ws1.Range("A16").Select() | Out-Null
